# DG: architecture pic update
# Remove the "Web" cloud shape and its dotted elbow connector from the
# architecture diagram on slide 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Delete the connector first (it references the cloud via endCxn),
# then delete the cloud shape itself.
$s.Shapes.Item("Elbow Connector 51").Delete()
$s.Shapes.Item("Cloud 50").Delete()
